$wb = $excel.ActiveWorkbook
$wsSettings = $wb.Worksheets.Item("Settings")
$ws = $wb.Worksheets.Item("Constants")

# Row 18: add Description for BrowserUrl
$ws.Range("C18").Value = "Url of the botsdna browser"

# Row 19: new TimeOut config entry
$ws.Range("A19").Value = "TimeOut"
$ws.Range("B19").Value = 5
$ws.Range("C19").Value = "Delay time for the check app state activity"

# Row 20: SystemException entry
$ws.Range("A20").Value = "SystemException"
$ws.Range("B20").Value = "Browser not loading"
$ws.Range("C20").Value = "Exception message to display in case of system exception"

# Row 21: BusinessException entry
$ws.Range("A21").Value = "BusinessException"
$ws.Range("B21").Value = "Server not found"
$ws.Range("C21").Value = "Exception message to display in case of Business exception"

# Row 22: OutputSheet entry
$ws.Range("A22").Value = "OutputSheet"
$ws.Range("B22").Value = "Sheet3"
$ws.Range("C22").Value = "Name of the sheet where the output excel will be written"

# Update the active selection to B23 as left by the author after data entry,
# then return focus to the Settings sheet (which remains the active tab).
$ws.Activate()
$ws.Range("B23").Select()
$wsSettings.Activate()
